# The "for"/"endfor" demonstration paragraph originally represents the
# M2Doc tags as real Word fields:
#   { fldChar begin } m: for v | self.  { fldChar end }
#   <bold red error message run>
#   { fldChar begin }  m: endfor  { fldChar end }
#
# TokenIteratorFieldRewriterSplit now emits those tags as plain literal text
# (brace-delimited) instead of Word fields, i.e.:
#   {m:for v | self.}<bold red error message run, prefixed with "    <---">{m:endfor}
#
# Find the paragraph that still uses the old field-based representation
# (it is the one whose range contains exactly the two fields) and replace
# its whole content with the new, field-free run layout.

$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Fields.Count -ge 2) {
        $target = $candidate
        break
    }
}
if ($target -eq $null) {
    # Fallback: the second paragraph in the known fixture.
    $target = $d.Paragraphs.Item(2)
}

$paragraphXml = @"
<w:p xmlns:w="$wNs" w:rsidR="007A2DC4" w:rsidRDefault="00B31BB7">
  <w:r><w:t>{m:</w:t></w:r>
  <w:r><w:t xml:space="preserve">for v </w:t></w:r>
  <w:r><w:t>|</w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t xml:space="preserve">self.}</w:t></w:r>
  <w:r>
    <w:rPr>
      <w:b w:val="on"/>
      <w:color w:val="FF0000"/>
    </w:rPr>
    <w:t xml:space="preserve">    &lt;---Invalid for statement: Expression "self." is invalid: missing feature access or service call</w:t>
  </w:r>
  <w:r><w:t>{</w:t></w:r>
  <w:r><w:t>m:</w:t></w:r>
  <w:r><w:t xml:space="preserve">endfor}</w:t></w:r>
</w:p>
"@

[void]$target.Range.InsertXML($paragraphXml)
